$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A8").Value = "Stretchy Taxi"
$ws.Range("B8").Value = "com.singleton.strechy"

[void]$ws.Range("B8").Select()
